$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MOSFET part name and datasheet/vendor link (row 6)
$ws.Range("A6").Value = "On Semiconductor 2N7000"
$ws.Range("H6").Value = "https://www.digikey.com/product-detail/en/on-semiconductor/2N7000/2N7000FS-ND/244278"

# Update the unit price for the MOSFET; the dependent formula in F6 (and the
# rolled-up total in F15) will recalculate automatically.
$ws.Range("C6").Value = 0.32

# Force recalculation so cached formula results stay in sync.
$excel.CalculateFull()

# Update the last active cell selection to match the saved view state.
$ws.Range("D16").Select()

$wb.Save()
